$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" quarterly sheet.
#    Seed it by copying the layout/styles of the existing "2022-Q2"
#    sheet (same headers, same column styling), then overwrite the
#    data cells with the new quarter's figures.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"

$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Range("A1:H3").Copy($newSheet.Range("A1:H3"))

# Row 2: 景顺长城中证500行业中性低波动指数 (003318)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'003318"
$newSheet.Range("C2").Value = "景顺长城中证500行业中性低波动指数"
$newSheet.Range("D2").Value = "'10.25"
$newSheet.Range("E2").Value = "'93.67"
$newSheet.Range("F2").Value = "'1.32"
$newSheet.Range("G2").Value = "'0.1353"
$newSheet.Range("H2").Value = 1

# Row 3: 华安中证500行业中性低波动ETF (512260)
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'512260"
$newSheet.Range("C3").Value = "华安中证500行业中性低波动ETF"
$newSheet.Range("D3").Value = "'1.07"
$newSheet.Range("E3").Value = "'97.91"
$newSheet.Range("F3").Value = "'1.39"
$newSheet.Range("G3").Value = "'0.0149"
$newSheet.Range("H3").Value = 1

# Position the new sheet right after "总计" and before "2022-Q2".
$newSheet.Move($templateSheet)

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: a new top data row for 2022-Q3 is
#    inserted, pushing the other quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Create row 6 (2020-Q4) by copying the style of the row above it, then
# fill in the shifted-down values.
$summary.Range("A5").Copy($summary.Range("A6"))
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2020-Q4"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.26

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 0.22

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.16

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.17

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.15

# Keep the originally-active "2020-Q4" tab selected (it stays the same
# underlying sheet; only its physical position/file shifted).
$wb.Worksheets.Item("2020-Q4").Activate()
